{"js": "// \"rollback multi-line if statement\": re-add the paragraphs that spell out\n// the 3rd {{if}} block's condition/body/endif on separate lines (plus the\n// blank spacer paragraphs around them), appended after the existing\n// paragraphs at the end of the document body.\nconst body = context.document.body;\n\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"@if CreateDate > 2020-12-31\", Word.InsertLocation.end);\nbody.insertParagraph(\"CreateDate is not less than 2021\", Word.InsertLocation.end);\nbody.insertParagraph(\"asdasdasdasdasdasd\", Word.InsertLocation.end);\nbody.insertParagraph(\"@endif\", Word.InsertLocation.end);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# \"rollback multi-line if statement\": re-add the paragraphs that spell out\n# the 3rd {{if}} block's condition/body/endif on separate lines (plus the\n# blank spacer paragraphs around them), appended after the existing\n# paragraphs at the end of the document body.\n$d = $word.ActiveDocument\n\n$d.Paragraphs.Add() | Out-Null\n$d.Paragraphs.Add() | Out-Null\n\n$p3 = $d.Paragraphs.Add()\n$p3.Range.Text = \"@if CreateDate > 2020-12-31\"\n\n$p4 = $d.Paragraphs.Add()\n$p4.Range.Text = \"CreateDate is not less than 2021\"\n\n$p5 = $d.Paragraphs.Add()\n$p5.Range.Text = \"asdasdasdasdasdasd\"\n\n$p6 = $d.Paragraphs.Add()\n$p6.Range.Text = \"@endif\"\n\n$d.Paragraphs.Add() | Out-Null\n"}
